$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New data for the debtor-detail table (rows 16-28).
# Columns: B=Tipo Doc, C=N° Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @("CC", "33333269", "MARY TORRES RIPOLL", "1809", 31249, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1808", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1807", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1806", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1805", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1804", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1803", 29509, 737717),
    @("CC", "33333269", "MARY TORRES RIPOLL", "1801", 29509, 737717),
    @("CC", "23002667", "BLANCA ROSA TORRES MUNOZ", "1809", 31249, 781242),
    @("CC", "23002667", "BLANCA ROSA TORRES MUNOZ", "1808", 31249, 781242),
    @("CC", "23002667", "BLANCA ROSA TORRES MUNOZ", "1807", 31249, 781242),
    @("CC", "23002667", "BLANCA ROSA TORRES MUNOZ", "1806", 31249, 781242),
    @("CC", "73157992", "SALVADOR FRIERI DEL CASTILLO", "1607", 80000, 2000000)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
